# Data update from DGS's 2021/08/13 report.
# Appends a new row (row 67) with the 2021/08/13 data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

# Column A holds a date-like label that is stored as literal text (matching
# every other cell in the column), displayed through the column's existing
# "yyyy/mm/dd" number format. Building it via a text formula and then
# pasting-as-values keeps it a plain string (not an auto-converted date
# serial) without creating a stray/unused number-format style.
$ws.Cells.Item($row, 1).NumberFormat = "yyyy/mm/dd"
$ws.Cells.Item($row, 1).Formula = '="2021/08/13"'
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).NumberFormat = "0.00"
$ws.Cells.Item($row, 2).Value = 319.9

$ws.Cells.Item($row, 3).NumberFormat = "0.00"
$ws.Cells.Item($row, 3).Value = 324.6

$ws.Cells.Item($row, 4).NumberFormat = "0.00"
$ws.Cells.Item($row, 4).Value = 0.95

$ws.Cells.Item($row, 5).NumberFormat = "0.00"
$ws.Cells.Item($row, 5).Value = 0.95
